$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.11671699999999
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 297.8183156666666
$ws.Range("N2").Value = 893.4549469999999
$ws.Range("O2").Value = 0.8852156413092672
$ws.Range("P2").Value = 0.8852156413092673
$ws.Range("Q2").Value = 8052.681343116555
$ws.Range("R2").Value = 72474.13208804899
$ws.Range("S2").Value = 0.06281656805163208
$ws.Range("T2").Value = 0.0628165680516321
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.11671699999999
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.07234764413494278
$ws.Range("P3").Value = 0.0723476441349428
$ws.Range("Q3").Value = 658.1362743231846
$ws.Range("R3").Value = 5923.226468908661
$ws.Range("S3").Value = 0.005133924999852255
$ws.Range("T3").Value = 0.005133924999852257
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.11671699999999
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.04243671455578994
$ws.Range("P4").Value = 0.04243671455578994
$ws.Range("Q4").Value = 386.040783306926
$ws.Range("R4").Value = 3474.367049762334
$ws.Range("S4").Value = 0.003011389138852937
$ws.Range("T4").Value = 0.003011389138852938
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 297.8183156666666
$ws.Range("N5").Value = 893.4549469999999
$ws.Range("O5").Value = 0.8852156413092672
$ws.Range("P5").Value = 0.8852156413092673
$ws.Range("Q5").Value = 102916.0565084721
$ws.Range("R5").Value = 926244.5085762488
$ws.Range("S5").Value = 0.8028174954167535
$ws.Range("T5").Value = 0.8028174954167536
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.07234764413494278
$ws.Range("P6").Value = 0.0723476441349428
$ws.Range("Q6").Value = 8411.209522952036
$ws.Range("R6").Value = 75700.88570656831
$ws.Range("S6").Value = 0.06561333956753401
$ws.Range("T6").Value = 0.06561333956753403
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.04243671455578994
$ws.Range("P7").Value = 0.04243671455578994
$ws.Range("Q7").Value = 4933.734910962486
$ws.Range("R7").Value = 44403.61419866237
$ws.Range("S7").Value = 0.03848659615074777
$ws.Range("T7").Value = 0.03848659615074778
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 297.8183156666666
$ws.Range("N8").Value = 893.4549469999999
$ws.Range("O8").Value = 0.8852156413092672
$ws.Range("P8").Value = 0.8852156413092673
$ws.Range("Q8").Value = 2510.232752901173
$ws.Range("R8").Value = 22592.09477611055
$ws.Range("S8").Value = 0.01958157784088166
$ws.Range("T8").Value = 0.01958157784088166
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.34034433333333
$ws.Range("N9").Value = 73.021033
$ws.Range("O9").Value = 0.07234764413494278
$ws.Range("P9").Value = 0.0723476441349428
$ws.Range("Q9").Value = 205.1584014423476
$ws.Range("R9").Value = 1846.425612981128
$ws.Range("S9").Value = 0.001600379567556514
$ws.Range("T9").Value = 0.001600379567556515
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.277234
$ws.Range("N10").Value = 42.831702
$ws.Range("O10").Value = 0.04243671455578994
$ws.Range("P10").Value = 0.04243671455578994
$ws.Range("Q10").Value = 120.339074268848
$ws.Range("R10").Value = 1083.051668419632
$ws.Range("S10").Value = 0.000938729266189229
$ws.Range("T10").Value = 0.0009387292661892295
